$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 updates
$ws.Range("M2").Value = 0.04009133333333333
$ws.Range("N2").Value = 0.120274
$ws.Range("O2").Value = 0.01033409631432067
$ws.Range("P2").Value = 0.01033409631432067
$ws.Range("Q2").Value = 0.001218936898666667
$ws.Range("R2").Value = 0.010970432088
$ws.Range("S2").Value = 0.01033409631432067
$ws.Range("T2").Value = 0.01033409631432067

# Row 3 updates
$ws.Range("O3").Value = 0.578569084147867
$ws.Range("P3").Value = 0.578569084147867
$ws.Range("S3").Value = 0.578569084147867
$ws.Range("T3").Value = 0.578569084147867

# Row 4 updates
$ws.Range("O4").Value = 0.4110968195378122
$ws.Range("P4").Value = 0.4110968195378122
$ws.Range("S4").Value = 0.4110968195378122
$ws.Range("T4").Value = 0.4110968195378122
